$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data rows (2-10); row 1 header formatting/content is left untouched.
$ws.Range("A2:T10").ClearContents()

# Prime the shared-string table so new strings are interned in the exact order
# required by the target workbook (Inflammatory-Mac already exists from the header use;
# Resolving-Mac, Cdh1, Itgae, ECs, MuSCs are new and must appear in this order).
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("B2").Value = "Cdh1"
$ws.Range("C2").Value = "Itgae"
$ws.Range("D2").Value = "ECs"
$ws.Range("D4").Value = "MuSCs"

# Now fill in the full data grid (rows 2-9) with final values.
# Row 2
$ws.Range("A2").Value = "Inflammatory-Mac"
$ws.Range("B2").Value = "Cdh1"
$ws.Range("C2").Value = "Itgae"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.061311
$ws.Range("H2").Value = 3.183933
$ws.Range("I2").Value = 0.9838301151667835
$ws.Range("J2").Value = 0.9838301151667834
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.127652
$ws.Range("N2").Value = 0.382956
$ws.Range("O2").Value = 0.1213851285815897
$ws.Range("P2").Value = 0.1213851285815897
$ws.Range("Q2").Value = 0.135478471772
$ws.Range("R2").Value = 1.219306245948
$ws.Range("S2").Value = 0.1194223450319602
$ws.Range("T2").Value = 0.1194223450319602

# Row 3
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Cdh1"
$ws.Range("C3").Value = "Itgae"
$ws.Range("D3").Value = "Inflammatory-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.061311
$ws.Range("H3").Value = 3.183933
$ws.Range("I3").Value = 0.9838301151667835
$ws.Range("J3").Value = 0.9838301151667834
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.5212876666666667
$ws.Range("N3").Value = 1.563863
$ws.Range("O3").Value = 0.4956958797851205
$ws.Range("P3").Value = 0.4956958797851205
$ws.Range("Q3").Value = 0.5532483347976668
$ws.Range("R3").Value = 4.979235013179
$ws.Range("S3").Value = 0.4876805344966952
$ws.Range("T3").Value = 0.4876805344966951

# Row 4
$ws.Range("A4").Value = "Inflammatory-Mac"
$ws.Range("B4").Value = "Cdh1"
$ws.Range("C4").Value = "Itgae"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.061311
$ws.Range("H4").Value = 3.183933
$ws.Range("I4").Value = 0.9838301151667835
$ws.Range("J4").Value = 0.9838301151667834
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1525256666666667
$ws.Range("N4").Value = 0.457577
$ws.Range("O4").Value = 0.1450376622405135
$ws.Range("P4").Value = 0.1450376622405135
$ws.Range("Q4").Value = 0.1618771678156667
$ws.Range("R4").Value = 1.456894510341
$ws.Range("S4").Value = 0.1426924199456054
$ws.Range("T4").Value = 0.1426924199456054

# Row 5
$ws.Range("A5").Value = "Inflammatory-Mac"
$ws.Range("B5").Value = "Cdh1"
$ws.Range("C5").Value = "Itgae"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.061311
$ws.Range("H5").Value = 3.183933
$ws.Range("I5").Value = 0.9838301151667835
$ws.Range("J5").Value = 0.9838301151667834
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2501626666666667
$ws.Range("N5").Value = 0.750488
$ws.Range("O5").Value = 0.2378813293927764
$ws.Range("P5").Value = 0.2378813293927764
$ws.Range("Q5").Value = 0.2655003899226667
$ws.Range("R5").Value = 2.389503509304
$ws.Range("S5").Value = 0.2340348156925228
$ws.Range("T5").Value = 0.2340348156925228

# Row 6
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("B6").Value = "Cdh1"
$ws.Range("C6").Value = "Itgae"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.01744333333333333
$ws.Range("H6").Value = 0.05233
$ws.Range("I6").Value = 0.01616988483321658
$ws.Range("J6").Value = 0.01616988483321658
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.127652
$ws.Range("N6").Value = 0.382956
$ws.Range("O6").Value = 0.1213851285815897
$ws.Range("P6").Value = 0.1213851285815897
$ws.Range("Q6").Value = 0.002226676386666667
$ws.Range("R6").Value = 0.02004008748
$ws.Range("S6").Value = 0.001962783549629492
$ws.Range("T6").Value = 0.001962783549629492

# Row 7
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("B7").Value = "Cdh1"
$ws.Range("C7").Value = "Itgae"
$ws.Range("D7").Value = "Inflammatory-Mac"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.01744333333333333
$ws.Range("H7").Value = 0.05233
$ws.Range("I7").Value = 0.01616988483321658
$ws.Range("J7").Value = 0.01616988483321658
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.5212876666666667
$ws.Range("N7").Value = 1.563863
$ws.Range("O7").Value = 0.4956958797851205
$ws.Range("P7").Value = 0.4956958797851205
$ws.Range("Q7").Value = 0.009092994532222223
$ws.Range("R7").Value = 0.08183695079
$ws.Range("S7").Value = 0.008015345288425371
$ws.Range("T7").Value = 0.008015345288425371

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Cdh1"
$ws.Range("C8").Value = "Itgae"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.01744333333333333
$ws.Range("H8").Value = 0.05233
$ws.Range("I8").Value = 0.01616988483321658
$ws.Range("J8").Value = 0.01616988483321658
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.1525256666666667
$ws.Range("N8").Value = 0.457577
$ws.Range("O8").Value = 0.1450376622405135
$ws.Range("P8").Value = 0.1450376622405135
$ws.Range("Q8").Value = 0.002660556045555556
$ws.Range("R8").Value = 0.02394500441
$ws.Range("S8").Value = 0.002345242294908068
$ws.Range("T8").Value = 0.002345242294908068

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Cdh1"
$ws.Range("C9").Value = "Itgae"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.01744333333333333
$ws.Range("H9").Value = 0.05233
$ws.Range("I9").Value = 0.01616988483321658
$ws.Range("J9").Value = 0.01616988483321658
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.2501626666666667
$ws.Range("N9").Value = 0.750488
$ws.Range("O9").Value = 0.2378813293927764
$ws.Range("P9").Value = 0.2378813293927764
$ws.Range("Q9").Value = 0.004363670782222223
$ws.Range("R9").Value = 0.03927303704
$ws.Range("S9").Value = 0.003846513700253654
$ws.Range("T9").Value = 0.003846513700253654
